# Fruta / hortaliza, semanal
#
# Inserts two new weekly price records for "Feria Lagunitas de Puerto Montt -
# Frutilla" at rows 118-119, pushing the previously-existing rows 118-122
# down to rows 120-124 (dimension grows from A1:T122 to A1:T124).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 118..122 down by two to make room for the new entries.
$ws.Rows.Item(118).Insert()
$ws.Rows.Item(118).Insert()

# New row 118: Primera quality, bandeja, Provincia de Melipilla.
$ws.Range("A118").Value = 4
$ws.Range("B118").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C118").Value = "Los Lagos"
$ws.Range("D118").Value = 44516
$ws.Range("E118").Value = 10
$ws.Range("F118").Value = "Fruta"
$ws.Range("G118").Value = 100101
$ws.Range("H118").Value = "Berries"
$ws.Range("I118").Value = 100112025
$ws.Range("J118").Value = "Frutilla"
$ws.Range("K118").Value = "Sin especificar"
$ws.Range("L118").Value = "Primera"
$ws.Range("M118").Value = 800
$ws.Range("N118").Value = 9000
$ws.Range("O118").Value = 9500
$ws.Range("P118").Value = 9250
$ws.Range("Q118").Value = "$/bandeja 7 kilos"
$ws.Range("R118").Value = "Provincia de Melipilla"
$ws.Range("S118").Value = 1321
$ws.Range("T118").Value = 7

# New row 119: Segunda quality, caja, Región de La Araucanía.
$ws.Range("A119").Value = 4
$ws.Range("B119").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C119").Value = "Los Lagos"
$ws.Range("D119").Value = 44516
$ws.Range("E119").Value = 10
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100101
$ws.Range("H119").Value = "Berries"
$ws.Range("I119").Value = 100112025
$ws.Range("J119").Value = "Frutilla"
$ws.Range("K119").Value = "Sin especificar"
$ws.Range("L119").Value = "Segunda"
$ws.Range("M119").Value = 800
$ws.Range("N119").Value = 8000
$ws.Range("O119").Value = 8500
$ws.Range("P119").Value = 8250
$ws.Range("Q119").Value = "$/caja 7 kilos"
$ws.Range("R119").Value = "Región de La Araucanía"
$ws.Range("S119").Value = 1179
$ws.Range("T119").Value = 7
